# Remove the trailing "Ver no Jupiter..." / copyright boilerplate that
# used to be scraped onto the end of the course page, along with the
# blank paragraph that separated it from the "Requisitos" text. The
# blank paragraph that precedes the final page-break paragraph is left
# in place.

$d = $word.ActiveDocument

# Locate the paragraph that holds "LOB1004: ..." (the last real content
# line of "Requisitos") via Find, then map the hit back to its paragraph
# index in the document's Paragraphs collection.
$hit = $d.Content
$found = $hit.Find.Execute("LOB1004", $false, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor text 'LOB1004' in the document"
}

$hitStart = $hit.Start
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($hitStart -ge $p.Range.Start -and $hitStart -lt $p.Range.End) {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not map the found text back to a paragraph"
}

# The three paragraphs right after the anchor are:
#   anchorIndex + 1 -> blank paragraph
#   anchorIndex + 2 -> "Ver no Jupiter Salvar em pdf Salvar em docx"
#   anchorIndex + 3 -> "© 2020 . Contact: ... Creative Commons Attribution"
# Delete that whole span in one go, keeping the blank paragraph that
# follows (right before the page-break paragraph).
$firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
$lastToRemove = $d.Paragraphs.Item($anchorIndex + 3)

$victim = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
$victim.Delete()

Write-Output "Removed boilerplate paragraphs after index $anchorIndex"
